# Refresh crypto price (D) / 1h-volume-change (E) snapshot cells to match the
# latest scrape, as produced by the scheduled GitHub Actions data-update job.
# All of these cells are stored as plain text in the sheet (prices use "."
# as a thousands separator in several rows, e.g. "29.613.00"), so every write
# below is forced to literal text -- a leading apostrophe stops Excel from
# "smart"-converting number-looking values (which would silently drop
# significant trailing zeros, e.g. "0.810" -> 0.81) -- and the style is reset
# right after so no stray quote-prefix/text number-format is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.613.00"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = "'1.600.29"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = "'212.32"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = "'26.89"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').Value = "'0.0601"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = "'1.830.04"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = "'1.602.02"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('D15').Value = "'29.622.29"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').Value = "'63.67"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = "'241.12"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('D19').Value = "'7.64"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.66%  '
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = "'9.27"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').Value = "'154.95"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('D31').Value = "'1.07"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = "'3.23"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = "'3.17"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('D34').Value = "'1.420.45"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('E36').Value = '  +4.53%  '
$ws.Range('E37').Value = '  -2.32%  '
$ws.Range('D38').Value = "'2.29"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  +3.19%  '
$ws.Range('D40').Value = "'0.544"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').Value = "'55.51"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.74%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = "'0.0494"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.09%  '
$ws.Range('D44').Value = "'0.810"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.52%  '
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = "'0.987"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +16.95%  '
$ws.Range('D47').Value = "'66.06"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D49').Value = "'1.740.29"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('D50').Value = "'86.16"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('E51').Value = '  +2.14%  '
